$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Settings sheet
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")

# The default HITL queue value is no longer pre-filled.
$settings.Range("B5").Value = ""

# Row 7 (OrchestratorAssetsFolder / Shared) is unchanged in content, row 8-13
# are brand-new rows describing the Action Center / HART task fields.
$settings.Range("A8").Value = "TaskTitle"
$settings.Range("B8").Value = "System Exception"
$settings.Range("C8").Value = "Title of the Action Center Task"

$settings.Range("A9").Value = "TaskCatalog"
$settings.Range("B9").Value = "SE"
$settings.Range("C9").Value = "Catalog of the Action Center Task"

$settings.Range("A10").Value = "Title"
$settings.Range("B10").Value = "Exception"
$settings.Range("C10").Value = "Title of the exception/case which is displayed in Action Center (used in the 'Create Form Task' activity)"

$settings.Range("A11").Value = "Description"
$settings.Range("B11").Value = "Detailed description of the exception"
$settings.Range("C11").Value = "Description of the exception/case which is displayed in Action Center (used in the 'Create Form Task' activity)"

$settings.Range("A12").Value = "SolutionDescription"
$settings.Range("B12").Value = "Description of the steps that should be performed in order to solve the issue"
$settings.Range("C12").Value = "Description of the steps performed to solve the exception/case which is displayed in Action Center (used in the 'Create Form Task' activity)"

$settings.Range("A13").Value = "FileName"
$settings.Range("B13").Value = "HART.pdf"
$settings.Range("C13").Value = "Name of the file displayed in Action Center (used in the 'Create Form Task' activity) The file should be stored in the Storage Bucket"

$settings.Activate()
$settings.Range("C7").Select()

# ---------------------------------------------------------------------------
# Assets sheet
# ---------------------------------------------------------------------------
$assets = $wb.Worksheets.Item("Assets")

$assets.Range("A3").Value = "TaskTitle"
$assets.Range("B3").Value = "HART_TaskTitle"
$assets.Range("C3").Value = "Shared"
$assets.Range("D3").Value = "Title of the Action Center Task"

$assets.Range("A4").Value = "TaskCatalog"
$assets.Range("B4").Value = "HART_TaskCatalog"
$assets.Range("C4").Value = "Shared"
$assets.Range("D4").Value = "Catalog of the Action Center Task"

$assets.Range("A5").Value = "Title"
$assets.Range("B5").Value = "HART_Title"
$assets.Range("C5").Value = "Shared"
$assets.Range("D5").Value = "Title of the exception/case which is displayed in Action Center (used in the 'Create Form Task' activity)"

$assets.Range("A6").Value = "Description"
$assets.Range("B6").Value = "HART_Description"
$assets.Range("C6").Value = "Shared"
$assets.Range("D6").Value = "Description of the exception/case which is displayed in Action Center (used in the 'Create Form Task' activity)"

$assets.Range("A7").Value = "SolutionDescription"
$assets.Range("B7").Value = "HART_SolutionDescription"
$assets.Range("C7").Value = "Shared"
$assets.Range("D7").Value = "Description of the steps performed to solve the exception/case which is displayed in Action Center (used in the 'Create Form Task' activity)"

$assets.Range("A8").Value = "FileName"
$assets.Range("B8").Value = "HART_FileName"
$assets.Range("C8").Value = "Shared"
$assets.Range("D8").Value = "Name of the file displayed in Action Center (used in the 'Create Form Task' activity) The file should be stored in the Storage Bucket"

# Two trailing placeholder rows near the bottom of the sheet were removed.
$assets.Rows.Item(993).Delete()
$assets.Rows.Item(993).Delete()

$assets.Activate()
$assets.Range("D13").Select()

# ---------------------------------------------------------------------------
# Re-select the Settings sheet as the active sheet/tab, matching the source
# workbook (tabSelected="1" on sheet1).
# ---------------------------------------------------------------------------
$settings.Activate()
